$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "25.762.03"
$ws.Range("E2").Value = "  -4.18%  "

$ws.Range("D3").Value = "1.813.27"
$ws.Range("E3").Value = "  -3.30%  "

$ws.Range("D4").Value = "'1.001"
$ws.Range("E4").Value = "  -0.15%  "

$ws.Range("D5").Value = "'275.94"
$ws.Range("E5").Value = "  -8.56%  "

$ws.Range("D6").Value = "'1.000"
$ws.Range("E6").Value = "  -0.14%  "

$ws.Range("D7").Value = "'0.5075"
$ws.Range("E7").Value = "  -4.68%  "

$ws.Range("D8").Value = "'0.3510"
$ws.Range("E8").Value = "  -6.54%  "

$ws.Range("E9").Value = "  -1.56%  "

$ws.Range("D10").Value = "'0.06678"
$ws.Range("E10").Value = "  -6.80%  "

$ws.Range("D11").Value = "'19.97"
$ws.Range("E11").Value = "  -7.52%  "

$ws.Range("D12").Value = "'0.8305"
$ws.Range("E12").Value = "  -6.24%  "

$ws.Range("D13").Value = "'0.07902"
$ws.Range("E13").Value = "  -2.88%  "

$ws.Range("D14").Value = "1.793.97"
$ws.Range("E14").Value = "  -4.25%  "

$ws.Range("D15").Value = "'5.071"
$ws.Range("E15").Value = "  -3.88%  "

$ws.Range("D16").Value = "'87.52"
$ws.Range("E16").Value = "  -6.00%  "

$ws.Range("D17").Value = "'1.000"
$ws.Range("E17").Value = "  -0.20%  "

$ws.Range("D18").Value = "'13.97"
$ws.Range("E18").Value = "  -5.32%  "

$ws.Range("D19").Value = "'0.000008039"
$ws.Range("E19").Value = "  -5.93%  "

$ws.Range("D20").Value = "'1.000"
$ws.Range("E20").Value = "  -0.13%  "

$ws.Range("D21").Value = "25.795.85"
$ws.Range("E21").Value = "  -4.76%  "

$ws.Range("D22").Value = "'4.717"
$ws.Range("E22").Value = "  -5.21%  "

$ws.Range("D23").Value = "'10.00"
$ws.Range("E23").Value = "  -6.47%  "

$ws.Range("D24").Value = "'6.037"
$ws.Range("E24").Value = "  -5.74%  "

$ws.Range("B25").Value = "Monero"
$ws.Range("C25").Value = "https://coinranking.com/coin/3mVx2FX_iJFp5+monero-xmr"
$ws.Range("D25").Value = "'140.99"
$ws.Range("E25").Value = "  -4.40%  "

$ws.Range("B26").Value = "LidoDAOToken"
$ws.Range("C26").Value = "https://coinranking.com/coin/Pe93bIOD2+lidodaotoken-ldo"
$ws.Range("D26").Value = "'2.187"
$ws.Range("E26").Value = "  -3.83%  "

$ws.Range("D27").Value = "'1.667"
$ws.Range("E27").Value = "  -4.40%  "

$ws.Range("E28").Value = "  -5.63%  "

$ws.Range("D29").Value = "'109.33"
$ws.Range("E29").Value = "  -4.66%  "

$ws.Range("E30").Value = "  -8.80%  "

$ws.Range("E31").Value = "  -8.05%  "

$ws.Range("D32").Value = "'0.08789"
$ws.Range("E32").Value = "  -3.43%  "

$ws.Range("D33").Value = "'0.04867"
$ws.Range("E33").Value = "  -2.48%  "

$ws.Range("B34").Value = "ImmutableX"
$ws.Range("C34").Value = "https://coinranking.com/coin/Z96jIvLU7+immutablex-imx"
$ws.Range("D34").Value = "'0.7258"
$ws.Range("E34").Value = "  -9.39%  "

$ws.Range("B35").Value = "ARBITRUM"
$ws.Range("C35").Value = "https://coinranking.com/coin/1Uo6s62Oc+arbitrum-arb"
$ws.Range("D35").Value = "'1.137"
$ws.Range("E35").Value = "  -3.09%  "

$ws.Range("D36").Value = "'2.879"
$ws.Range("E36").Value = "  -3.87%  "

$ws.Range("B37").Value = "MXToken"
$ws.Range("C37").Value = "https://coinranking.com/coin/QUC5kVAxSoB-+mxtoken-mx"
$ws.Range("D37").Value = "'3.157"
$ws.Range("E37").Value = "  -1.70%  "

$ws.Range("B38").Value = "Frax"
$ws.Range("C38").Value = "https://coinranking.com/coin/KfWtaeV1W+frax-frax"
$ws.Range("D38").Value = "'0.9997"
$ws.Range("E38").Value = "  -0.09%  "

$ws.Range("D39").Value = "'0.5180"
$ws.Range("E39").Value = "  -11.48%  "

$ws.Range("D40").Value = "'0.01838"
$ws.Range("E40").Value = "  -6.07%  "

$ws.Range("D41").Value = "'2.256"
$ws.Range("E41").Value = "  -12.96%  "

$ws.Range("D42").Value = "'0.9508"
$ws.Range("E42").Value = "  -11.09%  "

$ws.Range("D43").Value = "'112.85"
$ws.Range("E43").Value = "  -3.19%  "

$ws.Range("E44").Value = "  -6.85%  "

$ws.Range("D45").Value = "'8.054"
$ws.Range("E45").Value = "  -9.76%  "

$ws.Range("D46").Value = "'0.9997"
$ws.Range("E46").Value = "  -0.13%  "

$ws.Range("D47").Value = "'0.4551"
$ws.Range("E47").Value = "  -9.77%  "

$ws.Range("D48").Value = "'0.1356"
$ws.Range("E48").Value = "  -9.40%  "

$ws.Range("D49").Value = "'9.293"
$ws.Range("E49").Value = "  -6.62%  "

$ws.Range("D50").Value = "'36.20"
$ws.Range("E50").Value = "  -4.56%  "

$ws.Range("D51").Value = "'1.501"
$ws.Range("E51").Value = "  -6.81%  "
